$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 120 (pushes old rows 120..184 down to 121..185)
$ws.Rows.Item(120).Insert()

# Populate the newly inserted row 120 with the new weekly record
$ws.Range("A120").Value = 4
$ws.Range("B120").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C120").Value = "Los Lagos"
$ws.Range("D120").Value = 44488
$ws.Range("E120").Value = 10
$ws.Range("F120").Value = 100112045
$ws.Range("G120").Value = "Zapallo"
$ws.Range("H120").Value = "Paine"
$ws.Range("I120").Value = "1a (guarda)"
$ws.Range("J120").Value = 1100
$ws.Range("K120").Value = 400
$ws.Range("L120").Value = 450
$ws.Range("M120").Value = 425
$ws.Range("N120").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O120").Value = "Región Metropolitana"
$ws.Range("P120").Value = 425
$ws.Range("Q120").Value = 1
$ws.Range("R120").Value = "Hortaliza"

Write-Host ("UsedRange: " + $ws.UsedRange.Address())
